$wb = $excel.ActiveWorkbook

# Duplicate the "Bus_Makhulu_f" sheet to create the new multi-axle
# "Truck_Amandla_A1" sheet (same layout/conditional formatting as the
# other vehicle sheets), then update its Ackermann instance data.
$src = $wb.Worksheets.Item("Bus_Makhulu_f")
$src.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Truck_Amandla_A1"

# Instance name shown in H3
$new.Range("H3").Value = "Ackermann_Amandla_A1"

# Updated Rack values (x, y, z or scalar)
$new.Range("F6").Value = -0.9731
$new.Range("G6").Value = 0.55801
$new.Range("H6").Value = 2.5924

# Make the new sheet active / selected, matching where the author left off
[void]$new.Activate()
[void]$new.Range("J11").Select()
